$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H39").Value = 75.09090999999999
$ws.Range("I39").Value = 72.166664
$ws.Range("J39").Value = 78.59999999999999
$ws.Range("K39").Value = 216.499992
$ws.Range("L39").Value = 235.8
$ws.Range("M39").Value = 79.50000800000001
$ws.Range("N39").Value = -827.8
$ws.Range("H43").Value = 929.7778
$ws.Range("I43").Value = 892
$ws.Range("J43").Value = 940.5714
$ws.Range("K43").Value = 892
$ws.Range("L43").Value = 940.5714
$ws.Range("M43").Value = -823
$ws.Range("N43").Value = -1078.5714
$ws.Range("H76").Value = 4447444.5
$ws.Range("I76").Value = 5053414.5
$ws.Range("J76").Value = 3666.6667
$ws.Range("K76").Value = 5053414.5
$ws.Range("L76").Value = 3666.6667
$ws.Range("M76").Value = -5053099.5
$ws.Range("N76").Value = -4296.6667
$ws.Range("H79").Value = 4447444.5
$ws.Range("I79").Value = 5053414.5
$ws.Range("J79").Value = 3666.6667
$ws.Range("K79").Value = 5053414.5
$ws.Range("L79").Value = 3666.6667
$ws.Range("M79").Value = -5052322.5
$ws.Range("N79").Value = -5850.6667
$ws.Range("H125").Value = 7007612
$ws.Range("I125").Value = 507.75
$ws.Range("J125").Value = 9343313
$ws.Range("K125").Value = 4569.75
$ws.Range("L125").Value = 84089817
$ws.Range("M125").Value = -2109.75
$ws.Range("N125").Value = -84094737
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 30312644
$ws.Range("I97").Value = 33343808
$ws.Range("K97").Value = 33343808
$ws.Range("M97").Value = -33343312
$ws.Range("H102").Value = 2729
$ws.Range("I102").Value = 2688.7144
$ws.Range("J102").Value = 3011
$ws.Range("K102").Value = 2688.7144
$ws.Range("L102").Value = 3011
$ws.Range("M102").Value = -1066.7144
$ws.Range("N102").Value = -6255
$ws.Range("H109").Value = 200125.67
$ws.Range("J109").Value = 200125.67
$ws.Range("L109").Value = 200125.67
$ws.Range("N109").Value = -202899.67
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1640.2727
$ws.Range("I99").Value = 1297.8334
$ws.Range("J99").Value = 2051.2
$ws.Range("K99").Value = 1297.8334
$ws.Range("L99").Value = 2051.2
$ws.Range("M99").Value = 200.1666
$ws.Range("N99").Value = -5047.2
$ws.Range("H105").Value = 3554.4666
$ws.Range("I105").Value = 3530.6
$ws.Range("J105").Value = 3602.2
$ws.Range("K105").Value = 3530.6
$ws.Range("L105").Value = 3602.2
$ws.Range("M105").Value = -1783.6
$ws.Range("N105").Value = -7096.2
$ws.Range("H130").Value = 40833.332
$ws.Range("J130").Value = 40833.332
$ws.Range("L130").Value = 40833.332
$ws.Range("N130").Value = -50873.332
$ws.Range("H134").Value = 32261560
$ws.Range("I134").Value = 55557780
$ws.Range("J134").Value = 5255.3076
$ws.Range("K134").Value = 166673340
$ws.Range("L134").Value = 15765.9228
$ws.Range("M134").Value = -166670805
$ws.Range("N134").Value = -20835.9228
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 15683.277
$ws.Range("I62").Value = 22990
$ws.Range("K62").Value = 22990
$ws.Range("M62").Value = -22366
$ws.Range("H65").Value = 15683.277
$ws.Range("I65").Value = 22990
$ws.Range("K65").Value = 114950
$ws.Range("M65").Value = -111830
$ws.Range("H70").Value = 29990
$ws.Range("J70").Value = 29990
$ws.Range("L70").Value = 29990
$ws.Range("N70").Value = -30620
$ws.Range("H73").Value = 29990
$ws.Range("J73").Value = 29990
$ws.Range("L73").Value = 29990
$ws.Range("N73").Value = -32174
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 85026
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 113334.664
$ws.Range("K4").Value = 300
$ws.Range("L4").Value = 340003.992
$ws.Range("M4").Value = -188
$ws.Range("N4").Value = -340227.992
$ws.Range("H5").Value = 1118.5264
$ws.Range("I5").Value = 655.1429000000001
$ws.Range("J5").Value = 2416
$ws.Range("K5").Value = 1965.4287
$ws.Range("L5").Value = 7248
$ws.Range("M5").Value = -1853.4287
$ws.Range("N5").Value = -7472
$ws.Range("H23").Value = 515.96155
$ws.Range("I23").Value = 888.5833
$ws.Range("J23").Value = 196.57143
$ws.Range("K23").Value = 2665.7499
$ws.Range("L23").Value = 589.71429
$ws.Range("M23").Value = -2430.7499
$ws.Range("N23").Value = -1059.71429
$ws.Range("H68").Value = 15670
$ws.Range("J68").Value = 911.6
$ws.Range("L68").Value = 2734.8
$ws.Range("N68").Value = -4356.8
$ws.Range("H71").Value = 15670
$ws.Range("J71").Value = 911.6
$ws.Range("L71").Value = 8204.4
$ws.Range("N71").Value = -16316.4
$ws.Range("H107").Value = 311.58823
$ws.Range("J107").Value = 499.6
$ws.Range("L107").Value = 1498.8
$ws.Range("N107").Value = -5338.8
$ws.Range("H113").Value = 12821268
$ws.Range("I113").Value = 431
$ws.Range("J113").Value = 14706685
$ws.Range("K113").Value = 1293
$ws.Range("L113").Value = 44120055
$ws.Range("M113").Value = 877
$ws.Range("N113").Value = -44124395
$ws.Range("H122").Value = 952.8823
$ws.Range("J122").Value = 1071.3572
$ws.Range("L122").Value = 9642.2148
$ws.Range("N122").Value = -14542.2148
$ws.Range("H131").Value = 1456.2924
$ws.Range("J131").Value = 1547.2881
$ws.Range("L131").Value = 4641.8643
$ws.Range("N131").Value = -14721.8643
$ws.Range("H133").Value = 7684.5454
$ws.Range("I133").Value = 1632.5
$ws.Range("J133").Value = 11142.857
$ws.Range("K133").Value = 4897.5
$ws.Range("L133").Value = 33428.571
$ws.Range("M133").Value = 162.5
$ws.Range("N133").Value = -43548.571
$ws.Range("H135").Value = 1118.5264
$ws.Range("I135").Value = 655.1429000000001
$ws.Range("J135").Value = 2416
$ws.Range("K135").Value = 5896.2861
$ws.Range("L135").Value = 21744
$ws.Range("M135").Value = -3361.2861
$ws.Range("N135").Value = -26814
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 1000
$ws.Range("J5").Value = 1000
$ws.Range("L5").Value = 1000
$ws.Range("N5").Value = -1224
$ws.Range("H80").Value = 2699.8
$ws.Range("I80").Value = 2699.5
$ws.Range("J80").Value = 2700
$ws.Range("K80").Value = 2699.5
$ws.Range("L80").Value = 2700
$ws.Range("M80").Value = -1701.5
$ws.Range("N80").Value = -4696
$ws.Range("H83").Value = 2699.8
$ws.Range("I83").Value = 2699.5
$ws.Range("J83").Value = 2700
$ws.Range("K83").Value = 13497.5
$ws.Range("L83").Value = 13500
$ws.Range("M83").Value = -8505.5
$ws.Range("N83").Value = -23484
$ws.Range("H122").Value = 1112902.1
$ws.Range("I122").Value = 1853185
$ws.Range("J122").Value = 2477.75
$ws.Range("K122").Value = 5559555
$ws.Range("L122").Value = 7433.25
$ws.Range("M122").Value = -5557105
$ws.Range("N122").Value = -12333.25
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 70002
$ws.Range("J2").Value = 70002
$ws.Range("L2").Value = 70002
$ws.Range("N2").Value = -70226
$ws.Range("H7").Value = 3069.4092
$ws.Range("J7").Value = 3258.2354
$ws.Range("L7").Value = 3258.2354
$ws.Range("N7").Value = -3482.2354
$ws.Range("H22").Value = 10462
$ws.Range("I22").Value = 1196.6666
$ws.Range("J22").Value = 13936.5
$ws.Range("K22").Value = 1196.6666
$ws.Range("L22").Value = 13936.5
$ws.Range("M22").Value = -901.6666
$ws.Range("N22").Value = -14526.5
$ws.Range("H27").Value = 10462
$ws.Range("I27").Value = 1196.6666
$ws.Range("J27").Value = 13936.5
$ws.Range("K27").Value = 1196.6666
$ws.Range("L27").Value = 13936.5
$ws.Range("M27").Value = -1089.6666
$ws.Range("N27").Value = -14150.5
$ws.Range("H40").Value = 3473.842
$ws.Range("I40").Value = 3250.75
$ws.Range("J40").Value = 3533.3333
$ws.Range("K40").Value = 3250.75
$ws.Range("L40").Value = 3533.3333
$ws.Range("M40").Value = -3114.75
$ws.Range("N40").Value = -3805.3333
$ws.Range("H55").Value = 767.3333
$ws.Range("I55").Value = 800
$ws.Range("J55").Value = 751
$ws.Range("K55").Value = 800
$ws.Range("L55").Value = 751
$ws.Range("M55").Value = -627
$ws.Range("N55").Value = -1097
$ws.Range("H126").Value = 3069.4092
$ws.Range("J126").Value = 3258.2354
$ws.Range("L126").Value = 9774.706200000001
$ws.Range("N126").Value = -14714.7062
$ws.Range("H136").Value = 4934.3438
$ws.Range("I136").Value = 3187.35
$ws.Range("J136").Value = 7846
$ws.Range("K136").Value = 9562.049999999999
$ws.Range("L136").Value = 23538
$ws.Range("M136").Value = -7012.049999999999
$ws.Range("N136").Value = -28638
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 30781068
$ws.Range("J2").Value = 12984.25
$ws.Range("L2").Value = 12984.25
$ws.Range("N2").Value = -13208.25
$ws.Range("H17").Value = 7622
$ws.Range("I17").Value = 7622
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 7622
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -7450
$ws.Range("N17").ClearContents()
$ws.Range("H107").Value = 479.92856
$ws.Range("I107").Value = 453.66666
$ws.Range("K107").Value = 1360.99998
$ws.Range("M107").Value = 559.0000199999999
$ws.Range("H124").Value = 68429
$ws.Range("J124").Value = 68429
$ws.Range("L124").Value = 68429
$ws.Range("N124").Value = -78249
